$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$boundariesSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Update A2 (Version: ...)
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# Update A6 (Recommended Citation ...)
$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Tongde Coal Mine, China, M2174, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"
$aboutSheet.Range("A6").Value = $newCitation

# Update S2:S7 (build_version) on Boundaries and methane sources sheet
for ($row = 2; $row -le 7; $row++) {
    $boundariesSheet.Range("S" + $row).Value = $newVersion
}
